$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2833
$ws.Range("I62").Value = 2833
$ws.Range("K62").Value = 2833
$ws.Range("M62").Value = -2209
$ws.Range("H65").Value = 2833
$ws.Range("I65").Value = 2833
$ws.Range("K65").Value = 14165
$ws.Range("M65").Value = -11045
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H125").Value = 687.6
$ws.Range("I125").Value = 571
$ws.Range("J125").Value = 862.5
$ws.Range("K125").Value = 5139
$ws.Range("L125").Value = 7762.5
$ws.Range("M125").Value = -2679
$ws.Range("N125").Value = -12682.5
$ws.Range("H137").Value = 5299.6
$ws.Range("I137").Value = 4749.6665
$ws.Range("K137").Value = 14248.9995
$ws.Range("M137").Value = -11698.9995
$ws.Range("H138").Value = 3217.1853
$ws.Range("J138").Value = 3828.4
$ws.Range("L138").Value = 11485.2
$ws.Range("N138").Value = -21765.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2476.3333
$ws.Range("I61").Value = 2476.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2476.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2264.3333
$ws.Range("N61").Value = ""
$ws.Range("H74").Value = 2646.9524
$ws.Range("I74").Value = 2234.4707
$ws.Range("J74").Value = 4400
$ws.Range("K74").Value = 2234.4707
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -1360.4707
$ws.Range("N74").Value = -6148
$ws.Range("H77").Value = 2646.9524
$ws.Range("I77").Value = 2234.4707
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 11172.3535
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -6804.353499999999
$ws.Range("N77").Value = -30736
$ws.Range("H122").Value = 3859.8
$ws.Range("I122").Value = 4528.143
$ws.Range("J122").Value = 3275
$ws.Range("K122").Value = 13584.429
$ws.Range("L122").Value = 9825
$ws.Range("M122").Value = -11134.429
$ws.Range("N122").Value = -14725
$ws.Range("H136").Value = 2476.3333
$ws.Range("I136").Value = 2476.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7428.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4878.999899999999
$ws.Range("N136").Value = ""
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 675.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 675.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 675.5
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -1125.5
$ws.Range("H67").Value = 675.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 675.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 675.5
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -2235.5
$ws.Range("H86").Value = 1560
$ws.Range("I86").Value = 1366.6666
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1366.6666
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -243.6666
$ws.Range("N86").Value = -4096
$ws.Range("H89").Value = 1560
$ws.Range("I89").Value = 1366.6666
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 6833.333000000001
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -1217.333000000001
$ws.Range("N89").Value = -20482
$ws.Range("H134").Value = 8561.75
$ws.Range("I134").Value = 8248.5
$ws.Range("J134").Value = 8875
$ws.Range("K134").Value = 24745.5
$ws.Range("L134").Value = 26625
$ws.Range("M134").Value = -22210.5
$ws.Range("N134").Value = -31695
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 41.42857
$ws.Range("I19").Value = 41.42857
$ws.Range("K19").Value = 41.42857
$ws.Range("M19").Value = 128.57143
$ws.Range("H24").Value = 41.42857
$ws.Range("I24").Value = 41.42857
$ws.Range("K24").Value = 41.42857
$ws.Range("M24").Value = 128.57143
$ws.Range("H31").Value = 5782.4707
$ws.Range("I31").Value = 3505.5
$ws.Range("J31").Value = 7806.4443
$ws.Range("K31").Value = 3505.5
$ws.Range("L31").Value = 7806.4443
$ws.Range("M31").Value = -3210.5
$ws.Range("N31").Value = -8396.444299999999
$ws.Range("H34").Value = 5782.4707
$ws.Range("I34").Value = 3505.5
$ws.Range("J34").Value = 7806.4443
$ws.Range("K34").Value = 3505.5
$ws.Range("L34").Value = 7806.4443
$ws.Range("M34").Value = -3303.5
$ws.Range("N34").Value = -8210.444299999999
$ws.Range("H51").Value = 24142.857
$ws.Range("J51").Value = 24142.857
$ws.Range("L51").Value = 24142.857
$ws.Range("N51").Value = -25614.857
$ws.Range("H61").Value = 24142.857
$ws.Range("J61").Value = 24142.857
$ws.Range("L61").Value = 24142.857
$ws.Range("N61").Value = -24838.857
$ws.Range("H134").Value = 2847.3333
$ws.Range("I134").Value = 2999.2
$ws.Range("K134").Value = 8997.599999999999
$ws.Range("M134").Value = -6462.599999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 41
$ws.Range("I58").Value = 41
$ws.Range("K58").Value = 41
$ws.Range("M58").Value = 236
$ws.Range("H126").Value = 3346.8
$ws.Range("I126").Value = 3346.8
$ws.Range("K126").Value = 10040.4
$ws.Range("M126").Value = -7570.400000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 5667
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 5667
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -5531
$ws.Range("N40").Value = -1272
$ws.Range("H61").Value = 3857.5715
$ws.Range("I61").Value = 3801.2
$ws.Range("K61").Value = 3801.2
$ws.Range("M61").Value = -3599.2
$ws.Range("H113").Value = 3857.5715
$ws.Range("I113").Value = 3801.2
$ws.Range("K113").Value = 3801.2
$ws.Range("M113").Value = -1631.2
$ws.Range("H136").Value = 2391.2856
$ws.Range("I136").Value = 1559.75
$ws.Range("K136").Value = 4679.25
$ws.Range("M136").Value = -2129.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3271.4443
$ws.Range("I122").Value = 3777.8572
$ws.Range("J122").Value = 1499
$ws.Range("K122").Value = 11333.5716
$ws.Range("L122").Value = 4497
$ws.Range("M122").Value = -8883.571599999999
$ws.Range("N122").Value = -9397

Write-Output "Applied Kraken_Profits price/profit updates"
